$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (row 1) - force text format so date-like strings
# ("August 2024", "September 2024") are not auto-converted to date serials
$ws.Range("A1").NumberFormat = "@"
$ws.Range("A1").Value = "August 2024"
$ws.Range("G1").NumberFormat = "@"
$ws.Range("G1").Value = "September 2024"

# Update data values (row 2)
$ws.Range("A2").Value = 1.094
$ws.Range("B2").Value = -0.527
$ws.Range("C2").Value = -0.061
$ws.Range("D2").Value = -0.087
$ws.Range("E2").Value = -0.003
$ws.Range("F2").Value = -0.422
$ws.Range("G2").Value = -0.004
